# Applies the "commit after first publish" catalog update:
#  1) Updates the "Current Quantity" (column H) for a batch of existing
#     product rows.
#  2) Inserts a brand-new product row ("Urban Daily Kit") before the old
#     row 163 ("Vegetarian Miracle"), pushing that row and everything
#     below it down by one.
#
# All cells in this sheet are stored as text (every value - numbers,
# letters, SKUs - is a literal string), so every write below forces the
# cell's NumberFormat to "@" (Text) first. That keeps e.g. "37.00" or
# "-1" from being silently normalised into a numeric 37 / -1 by the
# usual Excel type-coercion that happens on a plain `.Value =` write.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
}

# ---------------------------------------------------------------------
# 1) Column H ("Current Quantity jBrookerStudio") updates.
# ---------------------------------------------------------------------
$hChanges = @(
    @{Row=17;  Value="5"},
    @{Row=19;  Value="9"},
    @{Row=25;  Value="3"},
    @{Row=30;  Value="4"},
    @{Row=31;  Value="3"},
    @{Row=32;  Value="5"},
    @{Row=34;  Value="5"},
    @{Row=35;  Value="3"},
    @{Row=40;  Value="4"},
    @{Row=51;  Value="7"},
    @{Row=59;  Value="4"},
    @{Row=61;  Value="3"},
    @{Row=62;  Value="3"},
    @{Row=63;  Value="2"},
    @{Row=67;  Value="3"},
    @{Row=70;  Value="11"},
    @{Row=72;  Value="4"},
    @{Row=74;  Value="2"},
    @{Row=76;  Value="2"},
    @{Row=79;  Value="5"},
    @{Row=80;  Value="4"},
    @{Row=82;  Value="4"},
    @{Row=83;  Value="4"},
    @{Row=89;  Value="3"},
    @{Row=90;  Value="4"},
    @{Row=93;  Value="3"},
    @{Row=95;  Value="5"},
    @{Row=96;  Value="4"},
    @{Row=97;  Value="6"},
    @{Row=99;  Value="2"},
    @{Row=102; Value="5"},
    @{Row=103; Value="5"},
    @{Row=110; Value="6"},
    @{Row=114; Value="1"},
    @{Row=116; Value="3"},
    @{Row=119; Value="3"},
    @{Row=124; Value="-1"},
    @{Row=125; Value="0"},
    @{Row=126; Value="3"},
    @{Row=127; Value="6"},
    @{Row=133; Value="12"},
    @{Row=134; Value="30"},
    @{Row=138; Value="8"},
    @{Row=146; Value="1"},
    @{Row=152; Value="5"},
    @{Row=153; Value="6"},
    @{Row=155; Value="9"},
    @{Row=156; Value="0"}
)

foreach ($chg in $hChanges) {
    Set-TextValue $ws.Cells.Item($chg.Row, 8) $chg.Value
}

# ---------------------------------------------------------------------
# 2) Insert the new "Urban Daily Kit" row at row 163 (old row 163 and
#    everything below it shifts down to row 164+ automatically).
# ---------------------------------------------------------------------
$ws.Rows.Item(163).Insert()

Set-TextValue $ws.Cells.Item(163, 1)  "ZJPNULDRSO3I3IGTJOZQSZIG"
Set-TextValue $ws.Cells.Item(163, 2)  "Urban Daily Kit"
Set-TextValue $ws.Cells.Item(163, 3)  ""
Set-TextValue $ws.Cells.Item(163, 4)  "Skin Regimen"
Set-TextValue $ws.Cells.Item(163, 5)  "DS-KIT"
Set-TextValue $ws.Cells.Item(163, 6)  "Regular"
Set-TextValue $ws.Cells.Item(163, 7)  "99.00"
Set-TextValue $ws.Cells.Item(163, 8)  ""
Set-TextValue $ws.Cells.Item(163, 9)  ""
Set-TextValue $ws.Cells.Item(163, 10) ""
Set-TextValue $ws.Cells.Item(163, 11) ""
Set-TextValue $ws.Cells.Item(163, 12) "Y"
